$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting existing rows 32:63 down to 33:64
$ws.Rows("32:32").Insert()

# Populate the new row 32 with the new data entry
$ws.Range("A32").Value2 = 5
$ws.Range("B32").Value2 = "Macroferia Regional de Talca"
$ws.Range("C32").Value2 = "Maule"
$ws.Range("D32").Value2 = 44557
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value2 = 7
$ws.Range("F32").Value2 = "Fruta"
$ws.Range("G32").Value2 = 100103
$ws.Range("H32").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I32").Value2 = 100103002
$ws.Range("J32").Value2 = "Ciruela"
$ws.Range("K32").Value2 = "Black Amber"
$ws.Range("L32").Value2 = "Primera"
$ws.Range("M32").Value2 = 350
$ws.Range("N32").Value2 = 17000
$ws.Range("O32").Value2 = 17000
$ws.Range("P32").Value2 = 17000
$ws.Range("Q32").Value2 = '$/bandeja 18 kilos granel'
$ws.Range("R32").Value2 = "Provincia de Curic" + [char]243
$ws.Range("S32").Value2 = 944
$ws.Range("T32").Value2 = 18
